$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" footer date text from 9/8/2015 to
#    10/21/2015 everywhere it appears: the slide master and every one of its
#    slide layouts.
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "9/8/2015") {
                $tr.Text = "10/21/2015"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Refresh the cover slide's title text box for the doc 1.2.1 / MSBR 1.0.2
#    release: "8th September, 2015" -> "21st November, 2015", version bumps.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

# Work right-to-left across the paragraph so earlier offsets stay valid.

# " v1.0.1" -> " v1.0.2"
$idx = $tr.Text.IndexOf(" v1.0.1")
if ($idx -ge 0) {
    $c = $tr.Characters($idx + 1, 7)
    $c.Text = " v1.0.2"
}

# ", 2015 -  doc version 1.2.0" -> ", 2015 -  doc version 1.2.1" (keeps the
# literal en-dash character already present in the run).
$idx = $tr.Text.IndexOf(", 2015")
if ($idx -ge 0) {
    $c = $tr.Characters($idx + 1, 27)
    $c.Text = ", 2015 " + [char]0x2013 + "  doc version 1.2.1"
}

# "September" -> "November"
$idx = $tr.Text.IndexOf("September")
if ($idx -ge 0) {
    $c = $tr.Characters($idx + 1, 9)
    $c.Text = "November"
}

# "8th " -> "21st "
$idx = $tr.Text.IndexOf("8th ")
if ($idx -ge 0) {
    $c = $tr.Characters($idx + 1, 4)
    $c.Text = "21st "
}
